# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Cebollín" (Vega Modelo de Temuco)
# at the top of the data block (new rows 208-209), pushing the existing
# records (previously rows 208-239) down to rows 210-241.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 208; this shifts rows
# 208:239 down to 210:241 and keeps cell formatting (e.g. the date style
# on column D) consistent with the row being pushed down.
$ws.Range("A208:A209").EntireRow.Insert()

# New row 208
$ws.Range("A208").Value = 10
$ws.Range("B208").Value = "Vega Modelo de Temuco"
$ws.Range("C208").Value = "La Araucanía"
$ws.Range("D208").Value = 44504
$ws.Range("E208").Value = 9
$ws.Range("F208").Value = 100112037
$ws.Range("G208").Value = "Cebollín"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 110
$ws.Range("K208").Value = 8000
$ws.Range("L208").Value = 8000
$ws.Range("M208").Value = 8000
$ws.Range("N208").Value = "`$/docena de paquetes"
$ws.Range("O208").Value = "Provincia de Cautín"
$ws.Range("P208").Value = 667
$ws.Range("Q208").Value = 12
$ws.Range("R208").Value = "Hortaliza"

# New row 209
$ws.Range("A209").Value = 10
$ws.Range("B209").Value = "Vega Modelo de Temuco"
$ws.Range("C209").Value = "La Araucanía"
$ws.Range("D209").Value = 44504
$ws.Range("E209").Value = 9
$ws.Range("F209").Value = 100112037
$ws.Range("G209").Value = "Cebollín"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 95
$ws.Range("K209").Value = 5000
$ws.Range("L209").Value = 5000
$ws.Range("M209").Value = 5000
$ws.Range("N209").Value = "`$/docena de paquetes"
$ws.Range("O209").Value = "Región de O'Higgins"
$ws.Range("P209").Value = 417
$ws.Range("Q209").Value = 12
$ws.Range("R209").Value = "Hortaliza"
